$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching formatting of existing headers (H1, etc.)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Data values for columns I (I0) and J (IF)
$values = @{
    2  = 8
    3  = 7
    4  = 7
    5  = 8
    6  = 8
    7  = 7
    8  = 8
    9  = 7
    10 = 6
    11 = 5
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
